$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) - columns B..E inputs, G is the sum.
$ws.Cells.Item(2, 2).Value = 3.230985683306322
$ws.Cells.Item(2, 3).Value = 1.667794583268128
$ws.Cells.Item(2, 4).Value = 337.1190423067083
$ws.Cells.Item(2, 5).Value = 645.3272768299601
$ws.Cells.Item(2, 7).Value = 987.3450994032428

$ws.Cells.Item(3, 2).Value = 0.04763786555579896
$ws.Cells.Item(3, 3).Value = 0.04240448674262143
$ws.Cells.Item(3, 4).Value = 26.21740644021617
$ws.Cells.Item(3, 5).Value = 645.3272768299601
$ws.Cells.Item(3, 7).Value = 671.6347256224747

$ws.Cells.Item(4, 2).Value = 1.459612070389937
$ws.Cells.Item(4, 3).Value = 0.3127903958511391
$ws.Cells.Item(4, 4).Value = 3.900430680208489
$ws.Cells.Item(4, 5).Value = 8.660232485948974
$ws.Cells.Item(4, 7).Value = 14.33306563239854

$ws.Cells.Item(5, 2).Value = 0.003994804209775715
$ws.Cells.Item(5, 3).Value = 10.29869402782916
$ws.Cells.Item(5, 4).Value = 689428.5527653177
$ws.Cells.Item(5, 5).Value = 616238.5361209477
$ws.Cells.Item(5, 7).Value = 1305677.391575098

$ws.Cells.Item(6, 2).Value = 1.459612070389937
$ws.Cells.Item(6, 3).Value = 1.667794583268128
$ws.Cells.Item(6, 4).Value = 26.21740644021617
$ws.Cells.Item(6, 5).Value = 8.660232485948974
$ws.Cells.Item(6, 7).Value = 38.00504557982321

$ws.Cells.Item(7, 2).Value = 3.230985683306322
$ws.Cells.Item(7, 3).Value = 1.667794583268128
$ws.Cells.Item(7, 4).Value = 0.8054896365839992
$ws.Cells.Item(7, 5).Value = 0.496779210170732
$ws.Cells.Item(7, 7).Value = 6.201049113329182

$ws.Cells.Item(8, 2).Value = [double]"2.074986032285508e-05"
$ws.Cells.Item(8, 3).Value = 0.002777888934908601
$ws.Cells.Item(8, 4).Value = 26.21740644021617
$ws.Cells.Item(8, 5).Value = 616238.5361209477
$ws.Cells.Item(8, 7).Value = 616264.7563260267

$ws.Cells.Item(9, 2).Value = 3.230985683306322
$ws.Cells.Item(9, 3).Value = 1.667794583268128
$ws.Cells.Item(9, 4).Value = 0.8054896365839992
$ws.Cells.Item(9, 5).Value = 645.3272768299601
$ws.Cells.Item(9, 7).Value = 651.0315467331185

$ws.Cells.Item(10, 2).Value = [double]"2.074986032285508e-05"
$ws.Cells.Item(10, 3).Value = 114.8270160096505
$ws.Cells.Item(10, 4).Value = 26.21740644021617
$ws.Cells.Item(10, 5).Value = 9353990175.932438
$ws.Cells.Item(10, 7).Value = 9353990316.976881
